# Apply "added placement_analysis data for my colleges" edit:
#  - Renumber the Sl.NO column (A8:A20) so the sequence is contiguous
#    (1..19 instead of 1..6,8..20), reflecting one fewer placement row.
#  - Move the active selection to H10 (from F6).
#  - Update the remembered sort range to match the new (smaller) extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber Sl.NO for rows 8 through 20 (each decreases by 1) ---
$ws.Range("A8").Value  = 7
$ws.Range("A9").Value  = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19

# --- Update the sheet's last-used sort range to reflect the new extent ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A23")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:A23"))
$ws.Sort.Apply()

# --- Move the active cell / selection to H10 ---
[void]$ws.Range("H10").Select()
